$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 62

# Use text format while assigning so Excel doesn't auto-convert the
# date/time-looking strings into date/time serial numbers, then restore
# the default "Normal" style so the cell carries no explicit style index.
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "2025-10-06"
$ws.Cells.Item($row, 1).Style = "Normal"

$ws.Cells.Item($row, 2).NumberFormat = "@"
$ws.Cells.Item($row, 2).Value = "21:21:51"
$ws.Cells.Item($row, 2).Style = "Normal"

$ws.Cells.Item($row, 3).NumberFormat = "@"
$ws.Cells.Item($row, 3).Value = "1.00 EUR = 1,778.7733"
$ws.Cells.Item($row, 3).Style = "Normal"
